$p = $ppt.ActivePresentation

# --- 1. Table style change on slide 5's table (graphicFrame) ---
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{27588D20-F093-45B7-97C3-7D590830AB2E}")

# --- 2. Swap the active theme's colour scheme from "Red Violet" (Integral)
#        back to the stock "Office" palette, mirroring the OOXML theme-part
#        swap (theme2.xml reverts to the plain Office Theme colours).
$tcs = $p.SlideMaster.Theme.ThemeColorScheme
$tcs.Item(1).RGB  = 0
$tcs.Item(2).RGB  = 16777215
$tcs.Item(3).RGB  = 6968388
$tcs.Item(4).RGB  = 15132391
$tcs.Item(5).RGB  = 13998939
$tcs.Item(6).RGB  = 3243501
$tcs.Item(7).RGB  = 10855845
$tcs.Item(8).RGB  = 49407
$tcs.Item(9).RGB  = 12874308
$tcs.Item(10).RGB = 4697456
$tcs.Item(11).RGB = 12673797
$tcs.Item(12).RGB = 7491477
